{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) in bold + color (#2C3E50) across specific resume bullet points,\n// matching the target diff exactly (only the 6 paragraphs below change;\n// visually-similar numbers elsewhere in the document, e.g. in the\n// PROFESSIONAL SUMMARY or KEY PROJECTS sections, must stay untouched).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Each entry: the paragraph's exact current text, plus the ordered list of\n// substrings within it that should become bold + colored. Order matters only\n// in that every term must be unique enough (within the paragraph) for\n// `search` to find exactly the intended occurrence(s).\nconst edits = [\n  {\n    text:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    terms: [\"23%\", \"64%\"],\n  },\n  {\n    text:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00B14.2% to \\u00B12.1%\",\n    terms: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"],\n  },\n  {\n    text: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    terms: [\"1,200\"],\n  },\n  {\n    text:\n      \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    terms: [\"$400M\", \"$1B\"],\n  },\n  {\n    text: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    terms: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    terms: [\"87%\", \"71%\"],\n  },\n];\n\nfor (const edit of edits) {\n  let paragraph = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === edit.text) {\n      paragraph = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!paragraph) {\n    continue;\n  }\n\n  for (const term of edit.terms) {\n    const found = paragraph.search(term, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < found.items.length; i++) {\n      found.items[i].font.bold = true;\n      found.items[i].font.color = \"#2C3E50\";\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) in bold + color (#2C3E50) across specific resume bullet points.\n# Only the 6 paragraphs enumerated below change; visually-similar numbers\n# elsewhere in the document (PROFESSIONAL SUMMARY, KEY PROJECTS \"Impact:\"\n# lines, etc.) must stay untouched, so each paragraph is located by its\n# exact original text before any term inside it is searched/formatted.\n\n$d = $word.ActiveDocument\n\n# Word's Font.Color is a BGR-ordered integer (0xBBGGRR), not RRGGBB -> the\n# target highlight color is #2C3E50.\n$highlightColor = 0x503E2C\n\n$bullet = [char]0x2022\n$plusMinus = [char]0x00B1\n\n$edits = @(\n    @{\n        Text  = \"$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Terms = @(\"23%\", \"64%\")\n    },\n    @{\n        Text  = \"$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ${plusMinus}4.2% to ${plusMinus}2.1%\"\n        Terms = @(\"87%\", \"71%\", \"${plusMinus}4.2%\", \"${plusMinus}2.1%\")\n    },\n    @{\n        Text  = \"$bullet Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Terms = @(\"1,200\")\n    },\n    @{\n        Text  = \"$bullet Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Terms = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Text  = \"$bullet Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Terms = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text  = \"$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Terms = @(\"87%\", \"71%\")\n    }\n)\n\nforeach ($edit in $edits) {\n    $targetParagraph = $null\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $edit.Text) {\n            $targetParagraph = $p\n            break\n        }\n    }\n    if ($null -eq $targetParagraph) {\n        continue\n    }\n\n    $paraStart = $targetParagraph.Range.Start\n    $paraEnd = $targetParagraph.Range.End\n\n    foreach ($term in $edit.Terms) {\n        $searchRange = $d.Range($paraStart, $paraEnd)\n        $found = $searchRange.Find.Execute($term)\n        if ($found) {\n            $searchRange.Bold = $true\n            $searchRange.Font.Color = $highlightColor\n        }\n    }\n}\n"}
